$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shield the touched cells from Excel's "looks like a number" auto-conversion
# while we write the new label text, then restore the default (unstyled) cell
# format so the output matches the original General-format cells exactly.
$targetRange = $ws.Range("B2,D2,E2,F2,G2,B3,C3,D3,E3,F3,D4,F4,G4,B6,C6,F6,G6,B7,E7,F7,G7,D8,F8")
$targetRange.NumberFormat = "@"

$ws.Range("B2").Value = "-0.05"
$ws.Range("D2").Value = "-0.32**"
$ws.Range("E2").Value = "-0.16"
$ws.Range("F2").Value = "-0.05"
$ws.Range("G2").Value = "-0.16"
$ws.Range("B3").Value = "-0.07"
$ws.Range("C3").Value = "-0.1"
$ws.Range("D3").Value = "-0.34***"
$ws.Range("E3").Value = "-0.13"
$ws.Range("F3").Value = "-0.05"
$ws.Range("D4").Value = "-0.16"
$ws.Range("F4").Value = "-0.16"
$ws.Range("G4").Value = "-0.32**"
$ws.Range("B6").Value = "-0.04"
$ws.Range("C6").Value = "-0.08"
$ws.Range("F6").Value = "-0.12"
$ws.Range("G6").Value = "-0.21"
$ws.Range("B7").Value = "-0.06"
$ws.Range("E7").Value = "-0.15"
$ws.Range("F7").Value = "-0.11"
$ws.Range("G7").Value = "-0.2"
$ws.Range("D8").Value = "-0.07"
$ws.Range("F8").Value = "-0.13"

$targetRange.Style = "Normal"
